$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# (e.g. "1.009", "18.10") are preserved exactly as in the source data
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "20.547.20"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "1.477.09"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "0.9750"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").Value = "279.21"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").Value = "0.3660"
$ws.Range("E7").Value = "  -1.90%  "
$ws.Range("D8").Value = "0.3074"
$ws.Range("E8").Value = "  -4.54%  "
$ws.Range("D9").Value = "39.68"
$ws.Range("E9").Value = "  -5.25%  "
$ws.Range("D10").Value = "1.060"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").Value = "0.06649"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "5.497"
$ws.Range("E13").Value = "  -3.54%  "
$ws.Range("D14").Value = "18.10"
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("D15").Value = "6.212"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").Value = "0.9747"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "0.00001032"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "1.475.10"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "0.05935"
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").Value = "69.58"
$ws.Range("E20").Value = "  -5.39%  "
$ws.Range("D21").Value = "5.462"
$ws.Range("E21").Value = "  -5.08%  "
$ws.Range("D22").Value = "14.53"
$ws.Range("E22").Value = "  -3.17%  "
$ws.Range("D23").Value = "11.05"
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("D24").Value = "2.243"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").Value = "20.593.59"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").Value = "140.98"
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").Value = "2.133"
$ws.Range("E27").Value = "  -9.40%  "
$ws.Range("D28").Value = "17.22"
$ws.Range("E28").Value = "  -3.17%  "
$ws.Range("D29").Value = "1.631.78"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").Value = "114.27"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").Value = "3.945"
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("D32").Value = "4.987"
$ws.Range("E32").Value = "  -8.07%  "
$ws.Range("D33").Value = "0.08031"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").Value = "0.8103"
$ws.Range("E34").Value = "  -5.25%  "
$ws.Range("D35").Value = "1.532"
$ws.Range("E35").Value = "  -6.34%  "
$ws.Range("D36").Value = "1.220"
$ws.Range("D37").Value = "0.05841"
$ws.Range("E37").Value = "  -4.91%  "
$ws.Range("D38").Value = "4.717"
$ws.Range("E38").Value = "  -5.60%  "
$ws.Range("D39").Value = "0.9743"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").Value = "0.02046"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("D41").Value = "7.649"
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("D42").Value = "10.40"
$ws.Range("E42").Value = "  -4.26%  "
$ws.Range("D43").Value = "0.1885"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").Value = "0.5296"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("D45").Value = "3.525"
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("E46").Value = "  -4.01%  "
$ws.Range("D47").Value = "119.24"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("D48").Value = "0.5202"
$ws.Range("E48").Value = "  -4.08%  "
$ws.Range("D49").Value = "1.795"
$ws.Range("E49").Value = "  -2.55%  "
$ws.Range("D50").Value = "0.06464"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").Value = "0.9976"
$ws.Range("E51").Value = "  -0.29%  "
